$d = $word.ActiveDocument

# --- 1. Bold the "Correct usage:" paragraph (pPr mark + run) ---
$d.Paragraphs(60).Range.Font.Bold = 1

# --- 2. Bold the "Wrong usage:" paragraph (pPr mark + run) ---
$d.Paragraphs(62).Range.Font.Bold = 1

# --- 3. Bold the "Definition:" run under "6. Interface" ---
$p55 = $d.Paragraphs(55).Range
$defStart = $p55.Start
$defEnd = $defStart + 11
$defR = $d.Range($defStart, $defEnd)
$defR.Font.Bold = 1

# --- 4. Bold the "Example:" run under "6. Interface" ---
$p56 = $d.Paragraphs(56).Range
$exStart = $p56.Start
$exEnd = $exStart + 8
$exR = $d.Range($exStart, $exEnd)
$exR.Font.Bold = 1

# --- 5. Move the _GoBack bookmark from the end of the document to wrap
#        the newly bolded "Definition:" run ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $d.Range($defStart, $defEnd))

# --- 6. Remove the "Day 1 end" text from the final paragraph, leaving it empty ---
$lastParaCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaCount).Range
$dayStart = $lastPara.Start
$dayEnd = $dayStart + 9
$dayR = $d.Range($dayStart, $dayEnd)
$dayR.Text = ""
